$d = $word.ActiveDocument

# Locate the unique sentence that needs to be rewritten.
$anchor = $d.Content
$anchor.Find.ClearFormatting()
$found = $anchor.Find.Execute("изменить номер кабинета, в котором работает Дерматолог", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not locate the target sentence to edit"
}

# $anchor now collapses to the matched range; expand it to the whole
# enclosing paragraph (minus the trailing paragraph mark) so we can
# rewrite the run layout in a single InsertXML call. Rewriting the
# whole paragraph body at once avoids relying on InsertXML placing
# content exactly at a mid-paragraph offset.
$para = $anchor.Paragraphs(1)
$paraStart = $para.Range.Start
$paraEnd = $para.Range.End - 1   # exclude the pilcrow
$full = $d.Range($paraStart, $paraEnd)

$rPr = '<w:rPr><w:rFonts w:ascii="Bookman Old Style" w:hAnsi="Bookman Old Style"/><w:color w:val="8EAADB" w:themeColor="accent1" w:themeTint="99"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="ru-RU"/></w:rPr>'

function Make-Run([string]$text, [bool]$preserve) {
    $escaped = $text.Replace("&", "&amp;").Replace("<", "&lt;").Replace(">", "&gt;")
    if ($preserve) {
        return "<w:r>$rPr<w:t xml:space=`"preserve`">$escaped</w:t></w:r>"
    } else {
        return "<w:r>$rPr<w:t>$escaped</w:t></w:r>"
    }
}

$runsXml = ""
$runsXml += Make-Run "изменить дату окончания карты на " $true
$runsXml += Make-Run "(05.11.2021) " $true
$runsXml += Make-Run "для всех старше определённой даты" $false
$runsXml += Make-Run " (" $true
$runsXml += Make-Run ">10.09.2020" $false
$runsXml += Make-Run ")" $false
$runsXml += Make-Run " начал" $true
$runsXml += Make-Run "а" $false
$runsXml += Make-Run " вед" $true
$runsXml += Make-Run "ения карты" $false

$xmlHeader = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p>'
$xmlFooter = '</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$full.InsertXML($xmlHeader + $runsXml + $xmlFooter)

Write-Host "Edit applied."
